# Updated test data for DC,TripCurrent, Voltdrop,BatteryStandby
#
# This script applies the content changes from the commit to both sheets:
#   - "Add Devices Loop A"      (sheet1)
#   - "Add Sounder Base Device" (sheet2)
#
# Changes per sheet:
#   F1 = "Loop"              G1 = "Column"
#   F2 = "Built-in Loop-A"   G2 = 2
#   F3 = "Built-in Loop-B"
#   F4 = "Built-in Loop-C"
#   F5 = "Built-in Loop-D"
#   B4 = "NGC-494/T412 OR TC-161" (replaces old "NGC-1191" text / blank cell),
#        with a border on left/right/bottom only (no top border)
#
# Plus the active selection on both sheets is moved to B4.

$wb = $excel.ActiveWorkbook

function Set-LoopColumns($ws) {
    $ws.Range("F1").Value = "Loop"
    $ws.Range("G1").Value = "Column"

    $ws.Range("F2").Value = "Built-in Loop-A"
    $ws.Range("G2").Value = 2

    $ws.Range("F3").Value = "Built-in Loop-B"
    $ws.Range("F4").Value = "Built-in Loop-C"
    $ws.Range("F5").Value = "Built-in Loop-D"

    # Match the existing header style used by row 7 (A7, B7, ...)
    $ws.Range("F1:G1").Style = $ws.Range("A7").Style
    # Match the existing bordered style already used in column F/G area (B4/F2:F5/G2)
    $ws.Range("F2:F5").Style = $ws.Range("F2").Style
    $ws.Range("G2").Style = $ws.Range("F2").Style
}

function Set-B4($ws) {
    $ws.Range("B4").Value = "NGC-494/T412 OR TC-161"

    $b4 = $ws.Range("B4")
    $b4.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $b4.Borders.Item(7).Weight = 2      # xlThin
    $b4.Borders.Item(10).LineStyle = 1  # xlEdgeRight
    $b4.Borders.Item(10).Weight = 2     # xlThin
    $b4.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
    $b4.Borders.Item(9).Weight = 2      # xlThin
    $b4.Borders.Item(8).LineStyle = 0   # xlEdgeTop -> none
}

$ws1 = $wb.Worksheets.Item("Add Devices Loop A")
Set-LoopColumns $ws1
Set-B4 $ws1
$ws1.Range("B4").Select()

$ws2 = $wb.Worksheets.Item("Add Sounder Base Device")
Set-LoopColumns $ws2
Set-B4 $ws2
$ws2.Range("B4").Select()

$ws1.Select()
